$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 16 (columns F:K), shifting cells F16:K44 up by one row.
# This removes the duplicate "extraTurn" attribute row that was left over
# from the playerColor/extraTurn pair, restoring the correct sequential
# sequence numbers in column G and re-aligning the shared formulas.
$ws.Range("F16:K16").Delete(-4162)

# Update the view so it matches the post-edit selection/scroll position.
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("G12:K12").Select()
